$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E3").Value = 24
$ws.Range("E14").Value = 35
$ws.Range("E16").Value = 300
$ws.Range("E18").Value = 91
